$wb = $excel.ActiveWorkbook

# Sheet1: MobileNumber, Date, Date&Time, Enquiry_Date, Enquiry_PhoneNumber,
#         User1RecId, Lead_PN, Sales_PN all changed to new generated values
#         (module/schedule re-run produced a fresh date/time + ids).
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws1.Range("F2").Value = "9428336107"
$ws1.Range("N2").Value = "2024-02-16"
$ws1.Range("P2").Value = "2024-02-16 03:47:47 PM"
$ws1.Range("AC2").Value = "2024-02-16"
$ws1.Range("AE2").Value = "2735621118"
$ws1.Range("AN2").Value = "94443"
$ws1.Range("AT2").Value = "6005539403"
$ws1.Range("AX2").Value = "7185939750"

# Sheet2: only MobileNumber, Enquiry_PhoneNumber, Lead_PN, Sales_PN changed.
$ws2 = $wb.Worksheets.Item("Sheet2")
$ws2.Range("F2").Value = "9428336107"
$ws2.Range("AE2").Value = "2735621118"
$ws2.Range("AT2").Value = "6005539403"
$ws2.Range("AX2").Value = "7185939750"

# Sheet3: only MobileNumber, Enquiry_PhoneNumber, Lead_PN, Sales_PN changed.
$ws3 = $wb.Worksheets.Item("Sheet3")
$ws3.Range("F2").Value = "9428336107"
$ws3.Range("AE2").Value = "2735621118"
$ws3.Range("AT2").Value = "6005539403"
$ws3.Range("AX2").Value = "7185939750"

# Sheet4: only MobileNumber, Enquiry_PhoneNumber, Lead_PN, Sales_PN changed.
$ws4 = $wb.Worksheets.Item("Sheet4")
$ws4.Range("F2").Value = "9428336107"
$ws4.Range("AE2").Value = "2735621118"
$ws4.Range("AT2").Value = "6005539403"
$ws4.Range("AX2").Value = "7185939750"
